$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attachments")
$ws.Range("F4").Value = "mapping between questions of the survey instrument and variables of the data set"
$ws.Range("F7").Value = "mapping between questions of the survey instrument and variables of the data set"
$ws.Range("F10").Value = "mapping between questions of the survey instrument and variables of the data set"
